{"js": "// The author replaced the literal phrase \"court order\" with\n// \"paternity or parentage order\" everywhere it appears in the body of the\n// Statement of Facts template (3 occurrences \u2014 the merge-field label near\n// \"by_court_order\", and twice in \"The court order {% if ...\" describing the\n// parenting-time order language).\nconst body = context.document.body;\n\nconst results = body.search(\"court order\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"paternity or parentage order\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The author replaced the literal phrase \"court order\" with\n# \"paternity or parentage order\" everywhere it appears in the body of the\n# Statement of Facts template (3 occurrences \u2014 the merge-field label near\n# \"by_court_order\", and twice in \"The court order {% if ...\" describing the\n# parenting-time order language).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Execute(\n    \"court order\",                    # FindText\n    $true,                             # MatchCase\n    $false,                            # MatchWholeWord\n    $false,                            # MatchWildcards\n    $false,                            # MatchSoundsLike\n    $false,                            # MatchAllWordForms\n    $true,                             # Forward\n    1,                                 # Wrap (wdFindContinue)\n    $false,                            # Format\n    \"paternity or parentage order\",    # ReplaceWith\n    2                                  # Replace (wdReplaceAll)\n)\n"}
